$wb = $excel.ActiveWorkbook

# --- india_wheat: selection change (O9 -> A3), drop topLeftCell scroll ---
$wsIndia = $wb.Worksheets.Item("india_wheat")
$wsIndia.Activate()
$wsIndia.Range("A3").Select()

# --- shifted_india_wheat_wo_anchors: fill in missing "COMMODITY:" label ---
$wsShiftedNoAnchor = $wb.Worksheets.Item("shifted_india_wheat_wo_anchors")
$wsShiftedNoAnchor.Activate()
$wsShiftedNoAnchor.Range("D6").Value = "COMMODITY:"
$wsShiftedNoAnchor.Range("E8").Select()

# --- e3: no longer the active tab (handled by later activations) ---
$wsE3 = $wb.Worksheets.Item("e3")
$wsE3.Activate()
$wsE3.Range("B3").Select()

# --- e3_shifted: selection change (I17 -> C17) ---
$wsE3Shifted = $wb.Worksheets.Item("e3_shifted")
$wsE3Shifted.Activate()
$wsE3Shifted.Range("C17").Select()

# --- Add e4 ---
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$e4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$e4.Name = "e4"

$e4.Range("C3").Value = "Main subject"
$e4.Range("C4").Value = "Heading 1"
$e4.Range("D4").Value = "Heading 2"
$e4.Range("E4").Value = "Heading 3"
$e4.Range("C5").Value = 12
$e4.Range("D5").Value = 2
$e4.Range("E5").Value = 1
$e4.Range("C6").Value = 23
$e4.Range("D6").Value = 4
$e4.Range("E6").Value = 2
$e4.Range("C7").Value = 34
$e4.Range("D7").Value = 6
$e4.Range("E7").Value = 3
$e4.Range("C8").Value = 45
$e4.Range("D8").Value = 8
$e4.Range("E8").Value = 4
$e4.Range("C9").Value = 56
$e4.Range("D9").Value = 10
$e4.Range("E9").Value = 5
$e4.Range("C10").Value = 67
$e4.Range("D10").Value = 10
$e4.Range("E10").Value = 6

$e4.Range("C11").Value = "Column 1"
$e4.Range("D11").Value = "Column 2"
$e4.Range("E11").Value = "Column 3"
$e4.Range("B11").Value = "anchor 1"

$e4.Range("C12").Value = 1
$e4.Range("D12").Value = 2
$e4.Range("E12").Value = 3
$e4.Range("C13").Value = 2
$e4.Range("D13").Value = 4
$e4.Range("E13").Value = 6
$e4.Range("C14").Value = 3
$e4.Range("D14").Value = 6
$e4.Range("E14").Value = 9
$e4.Range("C15").Value = 4
$e4.Range("D15").Value = 8
$e4.Range("E15").Value = 12
$e4.Range("C16").Value = 5
$e4.Range("D16").Value = 10
$e4.Range("E16").Value = 15

$e4.Range("B11").Select()

# --- Add e4_shifted (e4 shifted right by 1 column, down by 1 row) ---
$e4s = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $e4)
$e4s.Name = "e4_shifted"

$e4s.Range("D4").Value = "Main subject"
$e4s.Range("D5").Value = "Heading 1"
$e4s.Range("E5").Value = "Heading 2"
$e4s.Range("F5").Value = "Heading 3"
$e4s.Range("D6").Value = 12
$e4s.Range("E6").Value = 2
$e4s.Range("F6").Value = 1
$e4s.Range("D7").Value = 23
$e4s.Range("E7").Value = 4
$e4s.Range("F7").Value = 2
$e4s.Range("D8").Value = 34
$e4s.Range("E8").Value = 6
$e4s.Range("F8").Value = 3
$e4s.Range("D9").Value = 45
$e4s.Range("E9").Value = 8
$e4s.Range("F9").Value = 4
$e4s.Range("D10").Value = 56
$e4s.Range("E10").Value = 10
$e4s.Range("F10").Value = 5
$e4s.Range("D11").Value = 67
$e4s.Range("E11").Value = 10
$e4s.Range("F11").Value = 6

$e4s.Range("D12").Value = "Column 1"
$e4s.Range("E12").Value = "Column 2"
$e4s.Range("F12").Value = "Column 3"
$e4s.Range("C12").Value = "anchor 1"

$e4s.Range("D13").Value = 1
$e4s.Range("E13").Value = 2
$e4s.Range("F13").Value = 3
$e4s.Range("D14").Value = 2
$e4s.Range("E14").Value = 4
$e4s.Range("F14").Value = 6
$e4s.Range("D15").Value = 3
$e4s.Range("E15").Value = 6
$e4s.Range("F15").Value = 9
$e4s.Range("D16").Value = 4
$e4s.Range("E16").Value = 8
$e4s.Range("F16").Value = 12
$e4s.Range("D17").Value = 5
$e4s.Range("E17").Value = 10
$e4s.Range("F17").Value = 15

$e4s.Range("C12").Select()

# --- Add e4_misaligned (like e4_shifted, but 2nd table shifted further + extra column) ---
$e4m = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $e4s)
$e4m.Name = "e4_misaligned"

$e4m.Range("D4").Value = "Main subject"
$e4m.Range("D5").Value = "Heading 1"
$e4m.Range("E5").Value = "Heading 2"
$e4m.Range("F5").Value = "Heading 3"
$e4m.Range("D6").Value = 12
$e4m.Range("E6").Value = 2
$e4m.Range("F6").Value = 1
$e4m.Range("D7").Value = 23
$e4m.Range("E7").Value = 4
$e4m.Range("F7").Value = 2
$e4m.Range("D8").Value = 34
$e4m.Range("E8").Value = 6
$e4m.Range("F8").Value = 3
$e4m.Range("D9").Value = 45
$e4m.Range("E9").Value = 8
$e4m.Range("F9").Value = 4
$e4m.Range("D10").Value = 56
$e4m.Range("E10").Value = 10
$e4m.Range("F10").Value = 5
$e4m.Range("D11").Value = 67
$e4m.Range("E11").Value = 10
$e4m.Range("F11").Value = 6

$e4m.Range("F12").Value = "Column 1"
$e4m.Range("G12").Value = "Column 2"
$e4m.Range("H12").Value = "Column 3"
$e4m.Range("I12").Value = "Column 4"
$e4m.Range("C12").Value = "anchor 1"

$e4m.Range("F13").Value = 1
$e4m.Range("G13").Value = 2
$e4m.Range("H13").Value = 3
$e4m.Range("I13").Value = 4
$e4m.Range("F14").Value = 2
$e4m.Range("G14").Value = 4
$e4m.Range("H14").Value = 6
$e4m.Range("I14").Value = 8
$e4m.Range("F15").Value = 3
$e4m.Range("G15").Value = 6
$e4m.Range("H15").Value = 9
$e4m.Range("I15").Value = 12
$e4m.Range("F16").Value = 4
$e4m.Range("G16").Value = 8
$e4m.Range("H16").Value = 12
$e4m.Range("I16").Value = 16
$e4m.Range("F17").Value = 5
$e4m.Range("G17").Value = 10
$e4m.Range("H17").Value = 15
$e4m.Range("I17").Value = 20

$e4m.Activate()
$e4m.Range("F17").Select()
